# Generate Report for Handback
# Updates the handoff/handback/generate timestamps for the second file
# (75af15d5-9871-42fc-9627-421b87f9cd98) across the Overview, zh-cn and
# de-de sheets, reflecting a freshly regenerated handback report.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: "Latest HO Xliff Generate Date" for the 2nd file (row 3) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G3").Value = "2016-08-18 20:49:42"

# --- zh-cn sheet: "Correspond Handoff Datetime" (H) and
#     "Correspond Handback DateTime" (K) for the 2nd file (row 3) ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H3").Value = "2016-08-18 20:49:37"
$wsZhCn.Range("K3").Value = "2016-08-18 20:50:12"

# --- de-de sheet: "Correspond Handoff Datetime" (H) and
#     "Correspond Handback DateTime" (K) for the 2nd file (row 3) ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H3").Value = "2016-08-18 20:49:42"
$wsDeDe.Range("K3").Value = "2016-08-18 20:50:23"
